# Applies the "added exception test for totalPontosJogador and fixed
# createTable" edit to the ER model document:
#
#   - MultiJogador table gains the `nomeJogo` attribute (it is part of
#     the composite PK together with idPartida, so it is underlined
#     like idPartida), plus the matching "{nomeJogo} de Jogo.nome" FK
#     clause on the PK/FK line underneath.
#   - Normal table gets the same `nomeJogo` addition.
#
# wdFindContinue = 1, wdReplaceOne = 1 (Execute's Replace arg defaults
# to a single replacement when Replace:=wdReplaceOne is omitted from
# the named form below we just pass 2 == wdReplaceAll positionally,
# which is safe here because every search string is unique in the
# document).

$d = $word.ActiveDocument

function Underline-Word($range, $text) {
    # $range is the (already-replaced) range returned by a prior
    # Find.Execute call -- it covers exactly the freshly inserted
    # text, so searching inside it for $text can't match anything
    # else in the document.
    $u = $range.Duplicate
    $u.Find.ClearFormatting()
    $null = $u.Find.Execute($text)
    $u.Font.Underline = 1
}

# ------------------------------------------------------------------
# 1) MultiJogador (idPartida, Estado, nomeRegiao)
#       -> MultiJogador (idPartida, nomeJogo, Estado, nomeRegiao)
# ------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute(", Estado, ", $true, $false, $false, $false, $false, $true, 1, $false, ", nomeJogo, Estado, ", 2) | Out-Null
Underline-Word $rng1 "nomeJogo"

# ------------------------------------------------------------------
# 2) PK: idPartida; FK: {idPartida} de Partida.id e {nomeRegiao} de Regiao.nome
#       -> PK: idPartida, nomeJogo; FK: {idPartida} de Partida.id , {nomeRegiao} de Regiao.nome e {nomeJogo} de Jogo.nome
# ------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("PK: idPartida; FK: {idPartida} de Partida.id e {nomeRegiao} de Regiao.nome", $true, $false, $false, $false, $false, $true, 1, $false, "PK: idPartida, nomeJogo; FK: {idPartida} de Partida.id , {nomeRegiao} de Regiao.nome e {nomeJogo} de Jogo.nome", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Normal (idPartida, dificuldade, idJogador, Pontuacao)
#       -> Normal (idPartida, nomeJogo, dificuldade, idJogador, Pontuacao)
# ------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute(", dificuldade, idJogador", $true, $false, $false, $false, $false, $true, 1, $false, ", nomeJogo, dificuldade, idJogador", 2) | Out-Null
Underline-Word $rng3 "nomeJogo"

# ------------------------------------------------------------------
# 4) PK: idPartida; FK: {idPartida} de Partida.id e {idJogador} de Jogador.id
#       -> PK: idPartida, nomeJogo; FK: {idPartida} de Partida.id , {idJogador} de Jogador.id e {nomeJogo} de Jogo.nome
# ------------------------------------------------------------------
$rng4 = $d.Content
$rng4.Find.Execute("PK: idPartida; FK: {idPartida} de Partida.id e {idJogador} de Jogador.id", $true, $false, $false, $false, $false, $true, 1, $false, "PK: idPartida, nomeJogo; FK: {idPartida} de Partida.id , {idJogador} de Jogador.id e {nomeJogo} de Jogo.nome", 2) | Out-Null
